$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.524740333333333
$ws.Range("H2").Value = 4.574221
$ws.Range("I2").Value = 0.2062237893390968
$ws.Range("J2").Value = 0.2062237893390969
$ws.Range("O2").Value = 0.1304295628731768
$ws.Range("P2").Value = 0.1304295628731768
$ws.Range("Q2").Value = 1.797319687463667
$ws.Range("R2").Value = 16.175877187173
$ws.Range("S2").Value = 0.0268976786975485
$ws.Range("T2").Value = 0.02689767869754851
$ws.Range("G3").Value = 1.524740333333333
$ws.Range("H3").Value = 4.574221
$ws.Range("I3").Value = 0.2062237893390968
$ws.Range("J3").Value = 0.2062237893390969
$ws.Range("M3").Value = 4.539335
$ws.Range("N3").Value = 13.618005
$ws.Range("O3").Value = 0.5022718405737094
$ws.Range("P3").Value = 0.5022718405737095
$ws.Range("Q3").Value = 6.921307161011666
$ws.Range("R3").Value = 62.291764449105
$ws.Range("S3").Value = 0.1035804022414331
$ws.Range("T3").Value = 0.1035804022414331
$ws.Range("G4").Value = 1.524740333333333
$ws.Range("H4").Value = 4.574221
$ws.Range("I4").Value = 0.2062237893390968
$ws.Range("J4").Value = 0.2062237893390969
$ws.Range("M4").Value = 1.480144333333333
$ws.Range("N4").Value = 4.440433
$ws.Range("O4").Value = 0.1637761519293199
$ws.Range("P4").Value = 0.1637761519293199
$ws.Range("Q4").Value = 2.256835764188111
$ws.Range("R4").Value = 20.311521877693
$ws.Range("S4").Value = 0.03377453865423998
$ws.Range("T4").Value = 0.03377453865424
$ws.Range("G5").Value = 1.524740333333333
$ws.Range("H5").Value = 4.574221
$ws.Range("I5").Value = 0.2062237893390968
$ws.Range("J5").Value = 0.2062237893390969
$ws.Range("M5").Value = 1.839355666666667
$ws.Range("N5").Value = 5.518067
$ws.Range("O5").Value = 0.2035224446237938
$ws.Range("P5").Value = 0.2035224446237938
$ws.Range("Q5").Value = 2.804539772311889
$ws.Range("R5").Value = 25.240857950807
$ws.Range("S5").Value = 0.04197116974587525
$ws.Range("T5").Value = 0.04197116974587527
$ws.Range("I6").Value = 0.4308548451232278
$ws.Range("J6").Value = 0.4308548451232279
$ws.Range("O6").Value = 0.1304295628731768
$ws.Range("P6").Value = 0.1304295628731768
$ws.Range("S6").Value = 0.0561962091112129
$ws.Range("T6").Value = 0.05619620911121291
$ws.Range("I7").Value = 0.4308548451232278
$ws.Range("J7").Value = 0.4308548451232279
$ws.Range("M7").Value = 4.539335
$ws.Range("N7").Value = 13.618005
$ws.Range("O7").Value = 0.5022718405737094
$ws.Range("P7").Value = 0.5022718405737095
$ws.Range("Q7").Value = 14.460401171295
$ws.Range("R7").Value = 130.143610541655
$ws.Range("S7").Value = 0.2164062560801442
$ws.Range("T7").Value = 0.2164062560801442
$ws.Range("I8").Value = 0.4308548451232278
$ws.Range("J8").Value = 0.4308548451232279
$ws.Range("M8").Value = 1.480144333333333
$ws.Range("N8").Value = 4.440433
$ws.Range("O8").Value = 0.1637761519293199
$ws.Range("P8").Value = 0.1637761519293199
$ws.Range("Q8").Value = 4.715113744947
$ws.Range("R8").Value = 42.436023704523
$ws.Range("S8").Value = 0.07056374857438535
$ws.Range("T8").Value = 0.07056374857438538
$ws.Range("I9").Value = 0.4308548451232278
$ws.Range("J9").Value = 0.4308548451232279
$ws.Range("M9").Value = 1.839355666666667
$ws.Range("N9").Value = 5.518067
$ws.Range("O9").Value = 0.2035224446237938
$ws.Range("P9").Value = 0.2035224446237938
$ws.Range("Q9").Value = 5.859409106553001
$ws.Range("R9").Value = 52.734681958977
$ws.Range("S9").Value = 0.08768863135748539
$ws.Range("T9").Value = 0.08768863135748543
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2032796666666667
$ws.Range("H10").Value = 0.609839
$ws.Range("I10").Value = 0.02749392945088694
$ws.Range("J10").Value = 0.02749392945088694
$ws.Range("O10").Value = 0.1304295628731768
$ws.Range("P10").Value = 0.1304295628731768
$ws.Range("Q10").Value = 0.2396201759563333
$ws.Range("R10").Value = 2.156581583607
$ws.Range("S10").Value = 0.003586021199945145
$ws.Range("T10").Value = 0.003586021199945146
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.2032796666666667
$ws.Range("H11").Value = 0.609839
$ws.Range("I11").Value = 0.02749392945088694
$ws.Range("J11").Value = 0.02749392945088694
$ws.Range("M11").Value = 4.539335
$ws.Range("N11").Value = 13.618005
$ws.Range("O11").Value = 0.5022718405737094
$ws.Range("P11").Value = 0.5022718405737095
$ws.Range("Q11").Value = 0.9227545056883334
$ws.Range("R11").Value = 8.304790551195
$ws.Range("S11").Value = 0.0138094265499007
$ws.Range("T11").Value = 0.0138094265499007
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.2032796666666667
$ws.Range("H12").Value = 0.609839
$ws.Range("I12").Value = 0.02749392945088694
$ws.Range("J12").Value = 0.02749392945088694
$ws.Range("M12").Value = 1.480144333333333
$ws.Range("N12").Value = 4.440433
$ws.Range("O12").Value = 0.1637761519293199
$ws.Range("P12").Value = 0.1637761519293199
$ws.Range("Q12").Value = 0.3008832466985555
$ws.Range("R12").Value = 2.707949220287
$ws.Range("S12").Value = 0.004502849966882462
$ws.Range("T12").Value = 0.004502849966882464
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.2032796666666667
$ws.Range("H13").Value = 0.609839
$ws.Range("I13").Value = 0.02749392945088694
$ws.Range("J13").Value = 0.02749392945088694
$ws.Range("M13").Value = 1.839355666666667
$ws.Range("N13").Value = 5.518067
$ws.Range("O13").Value = 0.2035224446237938
$ws.Range("P13").Value = 0.2035224446237938
$ws.Range("Q13").Value = 0.3739036068014445
$ws.Range("R13").Value = 3.365132461213
$ws.Range("S13").Value = 0.005595631734158631
$ws.Range("T13").Value = 0.005595631734158633
$ws.Range("G14").Value = 2.061212666666667
$ws.Range("H14").Value = 6.183638
$ws.Range("I14").Value = 0.2787826080683977
$ws.Range("J14").Value = 0.2787826080683978
$ws.Range("O14").Value = 0.1304295628731768
$ws.Range("P14").Value = 0.1304295628731768
$ws.Range("Q14").Value = 2.429697716299333
$ws.Range("R14").Value = 21.867279446694
$ws.Range("S14").Value = 0.03636149370700528
$ws.Range("T14").Value = 0.0363614937070053
$ws.Range("G15").Value = 2.061212666666667
$ws.Range("H15").Value = 6.183638
$ws.Range("I15").Value = 0.2787826080683977
$ws.Range("J15").Value = 0.2787826080683978
$ws.Range("M15").Value = 4.539335
$ws.Range("N15").Value = 13.618005
$ws.Range("O15").Value = 0.5022718405737094
$ws.Range("P15").Value = 0.5022718405737095
$ws.Range("Q15").Value = 9.356534800243333
$ws.Range("R15").Value = 84.20881320219
$ws.Range("S15").Value = 0.1400246536744532
$ws.Range("T15").Value = 0.1400246536744532
$ws.Range("G16").Value = 2.061212666666667
$ws.Range("H16").Value = 6.183638
$ws.Range("I16").Value = 0.2787826080683977
$ws.Range("J16").Value = 0.2787826080683978
$ws.Range("M16").Value = 1.480144333333333
$ws.Range("N16").Value = 4.440433
$ws.Range("O16").Value = 0.1637761519293199
$ws.Range("P16").Value = 0.1637761519293199
$ws.Range("Q16").Value = 3.050892248361555
$ws.Range("R16").Value = 27.458030235254
$ws.Range("S16").Value = 0.04565794277426194
$ws.Range("T16").Value = 0.04565794277426197
$ws.Range("G17").Value = 2.061212666666667
$ws.Range("H17").Value = 6.183638
$ws.Range("I17").Value = 0.2787826080683977
$ws.Range("J17").Value = 0.2787826080683978
$ws.Range("M17").Value = 1.839355666666667
$ws.Range("N17").Value = 5.518067
$ws.Range("O17").Value = 0.2035224446237938
$ws.Range("P17").Value = 0.2035224446237938
$ws.Range("Q17").Value = 3.791303198638444
$ws.Range("R17").Value = 34.121728787746
$ws.Range("S17").Value = 0.05673851791267728
$ws.Range("T17").Value = 0.05673851791267731
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.4188103333333333
$ws.Range("H18").Value = 1.256431
$ws.Range("I18").Value = 0.05664482801839063
$ws.Range("J18").Value = 0.05664482801839064
$ws.Range("O18").Value = 0.1304295628731768
$ws.Range("P18").Value = 0.1304295628731768
$ws.Range("Q18").Value = 0.4936814754336667
$ws.Range("R18").Value = 4.443133278903
$ws.Range("S18").Value = 0.007388160157464968
$ws.Range("T18").Value = 0.007388160157464969
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.4188103333333333
$ws.Range("H19").Value = 1.256431
$ws.Range("I19").Value = 0.05664482801839063
$ws.Range("J19").Value = 0.05664482801839064
$ws.Range("M19").Value = 4.539335
$ws.Range("N19").Value = 13.618005
$ws.Range("O19").Value = 0.5022718405737094
$ws.Range("P19").Value = 0.5022718405737095
$ws.Range("Q19").Value = 1.901120404461667
$ws.Range("R19").Value = 17.110083640155
$ws.Range("S19").Value = 0.02845110202777829
$ws.Range("T19").Value = 0.0284511020277783
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.4188103333333333
$ws.Range("H20").Value = 1.256431
$ws.Range("I20").Value = 0.05664482801839063
$ws.Range("J20").Value = 0.05664482801839064
$ws.Range("M20").Value = 1.480144333333333
$ws.Range("N20").Value = 4.440433
$ws.Range("O20").Value = 0.1637761519293199
$ws.Range("P20").Value = 0.1637761519293199
$ws.Range("Q20").Value = 0.6198997416247778
$ws.Range("R20").Value = 5.579097674623
$ws.Range("S20").Value = 0.00927707195955014
$ws.Range("T20").Value = 0.009277071959550146
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.4188103333333333
$ws.Range("H21").Value = 1.256431
$ws.Range("I21").Value = 0.05664482801839063
$ws.Range("J21").Value = 0.05664482801839064
$ws.Range("M21").Value = 1.839355666666667
$ws.Range("N21").Value = 5.518067
$ws.Range("O21").Value = 0.2035224446237938
$ws.Range("P21").Value = 0.2035224446237938
$ws.Range("Q21").Value = 0.7703411598752223
$ws.Range("R21").Value = 6.933070438877
$ws.Range("S21").Value = 0.01152849387359723
$ws.Range("T21").Value = 0.01152849387359724
